# Update column F ("想去人数") values on the "展览" sheet and the
# corresponding rows on the "全部类型" sheet (which mirrors the same
# records one row further down).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row => new value, for the "展览" sheet
$exhibitUpdates = @{
    5  = 26
    8  = 13895
    10 = 82
    11 = 5602
    13 = 46
    17 = 68
    21 = 40
    22 = 10393
    23 = 1184
    24 = 25
    25 = 39
    26 = 3702
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row => new value, for the "全部类型" sheet (same records, shifted by 1 row)
$allUpdates = @{
    6  = 26
    9  = 13895
    11 = 82
    12 = 5602
    14 = 46
    18 = 68
    22 = 40
    24 = 10393
    25 = 1184
    26 = 25
    27 = 39
    28 = 3702
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
